# reports/european_funds_requests/templates/xlsx/template.xlsx
# fix: filter invalid items in eu funds basic packet
#
# The "Data" sheet's header row is reshuffled so that "Subscription Type"
# and "Technical Email" move earlier (right after the identifying columns),
# column A is widened, and the bestFit (auto best-fit) flag is dropped from
# columns A, D, E and F in favor of fixed custom widths. The active cell
# selection ends up on K1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Reorder header row (row 1) --------------------------------------------
# Before: Request ID | Created At      | Last Change At | Customer ID
#         | Customer Name | Customer Tax ID | Subscription Type
#         | Antivirus Quantity | EU Fund Packet | Technical Contact
#         | Technical Email | Postal Address
# After:  Request ID | Subscription Type | Created At | Customer ID
#         | Technical Email | Last Change At | Customer Name
#         | Customer Tax ID | Antivirus Quantity | EU Fund Packet
#         | Technical Contact | Postal Address
$headers = @(
    "Request ID",
    "Subscription Type",
    "Created At",
    "Customer ID",
    "Technical Email",
    "Last Change At",
    "Customer Name",
    "Customer Tax ID",
    "Antivirus Quantity",
    "EU Fund Packet",
    "Technical Contact",
    "Postal Address"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Column widths -----------------------------------------------------
# Widen column A and remove the bestFit auto-size marker from columns
# A, D, E and F by giving them an explicit custom width.
$ws.Columns.Item(1).ColumnWidth = 19           # column A: was auto best-fit, now a fixed, wider column
$ws.Columns.Item(4).ColumnWidth = 13.3333333   # column D: drop bestFit, keep same visual width
$ws.Columns.Item(5).ColumnWidth = 16.6666667   # column E: drop bestFit, keep same visual width
$ws.Columns.Item(6).ColumnWidth = 16.6666667   # column F: drop bestFit, keep same visual width

# --- Selection -----------------------------------------------------------
[void]$ws.Range("K1").Select()
